$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 173
$ws.Range("I4").Value = 173
$ws.Range("K4").Value = 173
$ws.Range("M4").Value = -59

$ws.Range("H33").Value = 757
$ws.Range("I33").Value = 528.73334
$ws.Range("K33").Value = 528.73334
$ws.Range("M33").Value = -299.73334

$ws.Range("H40").Value = 8559392
$ws.Range("I40").Value = 4622.625
$ws.Range("K40").Value = 4622.625
$ws.Range("M40").Value = -4447.625

$ws.Range("H76").Value = 7125.1
$ws.Range("I76").Value = 6036
$ws.Range("J76").Value = 9666.333000000001
$ws.Range("K76").Value = 6036
$ws.Range("L76").Value = 9666.333000000001
$ws.Range("M76").Value = -5721
$ws.Range("N76").Value = -10296.333

$ws.Range("H79").Value = 7125.1
$ws.Range("I79").Value = 6036
$ws.Range("J79").Value = 9666.333000000001
$ws.Range("K79").Value = 6036
$ws.Range("L79").Value = 9666.333000000001
$ws.Range("M79").Value = -4944
$ws.Range("N79").Value = -11850.333

$ws.Range("H112").Value = 36458.168
$ws.Range("I112").Value = 2493.5386
$ws.Range("K112").Value = 7480.6158
$ws.Range("M112").Value = -6372.6158

$ws.Range("H116").Value = 15322.777
$ws.Range("I116").Value = 15322.777
$ws.Range("K116").Value = 15322.777
$ws.Range("M116").Value = -11880.777

$ws.Range("H130").Value = 216666.5
$ws.Range("J130").Value = 216666.5
$ws.Range("L130").Value = 216666.5
$ws.Range("N130").Value = -226706.5

$ws.Range("H138").Value = 2877.1226
$ws.Range("I138").Value = 2281.6365
$ws.Range("J138").Value = 3362.3333
$ws.Range("K138").Value = 6844.9095
$ws.Range("L138").Value = 10086.9999
$ws.Range("M138").Value = -1704.9095
$ws.Range("N138").Value = -20366.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4191.1816
$ws.Range("I2").Value = 3274.25
$ws.Range("K2").Value = 3274.25
$ws.Range("M2").Value = -3161.25

$ws.Range("H5").Value = 132.6923
$ws.Range("I5").Value = 175.88889
$ws.Range("J5").Value = 35.5
$ws.Range("K5").Value = 175.88889
$ws.Range("L5").Value = 35.5
$ws.Range("M5").Value = -63.88889
$ws.Range("N5").Value = -259.5

$ws.Range("H31").Value = 2537.6
$ws.Range("I31").Value = 2537.6
$ws.Range("K31").Value = 2537.6
$ws.Range("M31").Value = -2243.6

$ws.Range("H61").Value = 27781748
$ws.Range("I61").Value = 28575398
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 28575398
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -28575186
$ws.Range("N61").Value = -4424

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null

$ws.Range("H116").Value = 4191.1816
$ws.Range("I116").Value = 3274.25
$ws.Range("K116").Value = 3274.25
$ws.Range("M116").Value = -980.25

$ws.Range("H122").Value = 1919.25
$ws.Range("I122").Value = 2184.6667
$ws.Range("J122").Value = 1123
$ws.Range("K122").Value = 6554.000100000001
$ws.Range("L122").Value = 3369
$ws.Range("M122").Value = -4104.000100000001
$ws.Range("N122").Value = -8269

$ws.Range("H136").Value = 27781748
$ws.Range("I136").Value = 28575398
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 85726194
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -85723644
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 60602
$ws.Range("J2").Value = 60602
$ws.Range("L2").Value = 60602
$ws.Range("N2").Value = -60828

$ws.Range("H3").Value = 4191.1816
$ws.Range("I3").Value = 3274.25
$ws.Range("K3").Value = 3274.25
$ws.Range("M3").Value = -3160.25

$ws.Range("H4").Value = 132.6923
$ws.Range("I4").Value = 175.88889
$ws.Range("J4").Value = 35.5
$ws.Range("K4").Value = 175.88889
$ws.Range("L4").Value = 35.5
$ws.Range("M4").Value = -60.88889
$ws.Range("N4").Value = -265.5

$ws.Range("H22").Value = 3664201.5
$ws.Range("I22").Value = 484.875
$ws.Range("K22").Value = 484.875
$ws.Range("M22").Value = -311.875

$ws.Range("H64").Value = 499.5
$ws.Range("I64").Value = 499.5
$ws.Range("K64").Value = 499.5
$ws.Range("M64").Value = -274.5

$ws.Range("H67").Value = 499.5
$ws.Range("I67").Value = 499.5
$ws.Range("K67").Value = 499.5
$ws.Range("M67").Value = 280.5

$ws.Range("H74").Value = 33329.332
$ws.Range("J74").Value = 33329.332
$ws.Range("L74").Value = 33329.332
$ws.Range("N74").Value = -35201.332

$ws.Range("H77").Value = 33329.332
$ws.Range("J77").Value = 33329.332
$ws.Range("L77").Value = 99987.99600000001
$ws.Range("N77").Value = -109347.996

$ws.Range("H94").Value = 649
$ws.Range("I94").Value = 524.5
$ws.Range("J94").Value = 773.5
$ws.Range("K94").Value = 524.5
$ws.Range("L94").Value = 773.5
$ws.Range("M94").Value = -73.5
$ws.Range("N94").Value = -1675.5

$ws.Range("H105").Value = 2356.2222
$ws.Range("I105").Value = 2025.875
$ws.Range("J105").Value = 4999
$ws.Range("K105").Value = 2025.875
$ws.Range("L105").Value = 4999
$ws.Range("M105").Value = -278.875
$ws.Range("N105").Value = -8493

$ws.Range("H134").Value = 10640713
$ws.Range("I134").Value = 11906522
$ws.Range("K134").Value = 35719566
$ws.Range("M134").Value = -35717031

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 16853
$ws.Range("I41").Value = 12723.6
$ws.Range("J41").Value = 37500
$ws.Range("K41").Value = 12723.6
$ws.Range("L41").Value = 37500
$ws.Range("M41").Value = -12295.6
$ws.Range("N41").Value = -38356

$ws.Range("H105").Value = 1819288.1
$ws.Range("I105").Value = 1819288.1
$ws.Range("K105").Value = 1819288.1
$ws.Range("M105").Value = -1817541.1

$ws.Range("H107").Value = 51848.7
$ws.Range("J107").Value = 94019.45
$ws.Range("L107").Value = 94019.45
$ws.Range("N107").Value = -97859.45

$ws.Range("H108").Value = 75000
$ws.Range("J108").Value = 75000
$ws.Range("L108").Value = 75000
$ws.Range("N108").Value = -82680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 328.0909
$ws.Range("J2").Value = 522.25
$ws.Range("L2").Value = 3133.5
$ws.Range("N2").Value = -3359.5

$ws.Range("H44").Value = 279.8
$ws.Range("I44").Value = 149.5
$ws.Range("K44").Value = 448.5
$ws.Range("M44").Value = -50.5

$ws.Range("H47").Value = 537.3333
$ws.Range("I47").Value = 537.3333
$ws.Range("K47").Value = 1611.9999
$ws.Range("M47").Value = -1180.9999

$ws.Range("H92").Value = 865
$ws.Range("J92").Value = 1500
$ws.Range("L92").Value = 4500
$ws.Range("N92").Value = -6996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 50000
$ws.Range("I18").Value = 50000
$ws.Range("K18").Value = 50000
$ws.Range("M18").Value = -49707

$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4999.6
$ws.Range("I22").Value = 5749.5
$ws.Range("J22").Value = 4499.6665
$ws.Range("K22").Value = 5749.5
$ws.Range("L22").Value = 4499.6665
$ws.Range("M22").Value = -5454.5
$ws.Range("N22").Value = -5089.6665

$ws.Range("H27").Value = 4999.6
$ws.Range("I27").Value = 5749.5
$ws.Range("J27").Value = 4499.6665
$ws.Range("K27").Value = 5749.5
$ws.Range("L27").Value = 4499.6665
$ws.Range("M27").Value = -5642.5
$ws.Range("N27").Value = -4713.6665

$ws.Range("H128").Value = 69999
$ws.Range("J128").Value = 69999
$ws.Range("L128").Value = 69999
$ws.Range("N128").Value = -79959

$ws.Range("H136").Value = 3091.1667
$ws.Range("I136").Value = 2516.6667
$ws.Range("K136").Value = 7550.000100000001
$ws.Range("M136").Value = -5000.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 7349.3335
$ws.Range("I55").Value = 2048
$ws.Range("K55").Value = 2048
$ws.Range("M55").Value = -1771

$ws.Range("H96").Value = 1697.32
$ws.Range("I96").Value = 1409.4615
$ws.Range("J96").Value = 2009.1666
$ws.Range("K96").Value = 1409.4615
$ws.Range("L96").Value = 2009.1666
$ws.Range("M96").Value = -36.46149999999989
$ws.Range("N96").Value = -4755.1666

$ws.Range("H105").Value = 21152
$ws.Range("J105").Value = 21728.5
$ws.Range("L105").Value = 21728.5
$ws.Range("N105").Value = -28716.5

$ws.Range("H141").Value = 79631.336
$ws.Range("J141").Value = 79631.336
$ws.Range("L141").Value = 79631.336
$ws.Range("N141").Value = -89991.336
